$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "51.599.77"
$ws.Range("E2").Value = "  -1.54%  "

# Row 3
$ws.Range("D3").Value = "2.948.34"
$ws.Range("E3").Value = "  -2.24%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "374.57"
$ws.Range("E5").Value = "  +5.12%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.31"
$ws.Range("E6").Value = "  -3.71%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.549"
$ws.Range("E7").Value = "  -2.98%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.597"
$ws.Range("E9").Value = "  -4.86%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.47"
$ws.Range("E10").Value = "  -3.17%  "

# Row 11
$ws.Range("E11").Value = "  -0.31%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0843"
$ws.Range("E12").Value = "  -2.32%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.50"
$ws.Range("E13").Value = "  -4.87%  "

# Row 14
$ws.Range("D14").Value = "3.408.48"
$ws.Range("E14").Value = "  -2.16%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.44"
$ws.Range("E15").Value = "  -4.75%  "

# Row 16
$ws.Range("D16").Value = "2.938.04"
$ws.Range("E16").Value = "  -2.35%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.945"
$ws.Range("E17").Value = "  -8.47%  "

# Row 18
$ws.Range("D18").Value = "51.577.62"
$ws.Range("E18").Value = "  -1.70%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.32"
$ws.Range("E19").Value = "  -7.08%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.35"
$ws.Range("E20").Value = "  -3.12%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.13"
$ws.Range("E21").Value = "  -4.52%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0952"
$ws.Range("E22").Value = "  -2.67%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.86"
$ws.Range("E23").Value = "  -1.27%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "262.40"
$ws.Range("E24").Value = "  -1.30%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.72"
$ws.Range("E25").Value = "  -1.39%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.173"
$ws.Range("E26").Value = "  -4.04%  "

# Row 27
$ws.Range("E27").Value = "  +0.00%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "26.01"
$ws.Range("E28").Value = "  -4.47%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.22"
$ws.Range("E29").Value = "  -6.29%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.85"
$ws.Range("E30").Value = "  +5.92%  "

# Row 31
$ws.Range("E31").Value = "  -3.61%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.96"
$ws.Range("E32").Value = "  -4.18%  "

# Row 33
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.98"
$ws.Range("E33").Value = "  -5.55%  "

# Row 34
$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.11"
$ws.Range("E34").Value = "  -3.71%  "

# Row 35
$ws.Range("E35").Value = "  +0.25%  "

# Row 36
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0430"
$ws.Range("E36").Value = "  -3.70%  "

# Row 37
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.49%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.04"
$ws.Range("E38").Value = "  -6.29%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.30"
$ws.Range("E39").Value = "  -4.89%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.60"
$ws.Range("E40").Value = "  -4.34%  "

# Row 41
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.87"
$ws.Range("E41").Value = "  -6.64%  "

# Row 42
$ws.Range("E42").Value = "  -3.36%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.40"
$ws.Range("E43").Value = "  -3.45%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.67"
$ws.Range("E44").Value = "  -3.40%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.10"
$ws.Range("E45").Value = "  -2.29%  "

# Row 46
$ws.Range("D46").Value = "2.037.33"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.22"
$ws.Range("E48").Value = "  -5.73%  "

# Row 49
$ws.Range("E49").Value = "  +4.64%  "

# Row 50
$ws.Range("D50").Value = "3.230.01"
$ws.Range("E50").Value = "  -2.10%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0326"
$ws.Range("E51").Value = "  -3.59%  "
